# Add four new Time Recording Log entries (rows 21-24) to the "작성자명"
# sheet. Existing rows 21-24 are blank template rows that already carry the
# correct cell styles, so only the values need to be written.
#
# The new shared-string table entries are appended in the exact order the
# author entered them: all four Date cells (column A, top to bottom) first,
# then the Activity cells (column F) for rows 24, 23 and 21 (row 22 reuses
# an Activity description already present earlier in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Date column (A21:A24) ----
$ws.Range("A21").Value = "10월 27일"
$c = $ws.Range("A21").Characters(4, 4)
$c.Font.Name = "돋움"
$c.Font.Size = 10
$c.Font.ColorIndex = -4105

$ws.Range("A22").Value = "10월 29일"
$c = $ws.Range("A22").Characters(4, 4)
$c.Font.Name = "돋움"
$c.Font.Size = 10
$c.Font.ColorIndex = -4105

$ws.Range("A23").Value = "11월 1일"
$c = $ws.Range("A23").Characters(4, 3)
$c.Font.Name = "돋움"
$c.Font.Size = 10
$c.Font.ColorIndex = -4105

$ws.Range("A24").Value = "11월 2일"
$c = $ws.Range("A24").Characters(4, 3)
$c.Font.Name = "돋움"
$c.Font.Size = 10
$c.Font.ColorIndex = -4105

# ---- Start / Stop / Interruption Time / Delta Time (B:E, rows 21-24) ----
# Row 21: 14:00-16:30, interrupt 30, delta 120
$ws.Range("B21").Value = 0.58333333333333337
$ws.Range("C21").Value = 0.6875
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = 120

# Row 22: 19:00-20:30, interrupt 0, delta 90
$ws.Range("B22").Value = 0.79166666666666663
$ws.Range("C22").Value = 0.85416666666666663
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 90

# Row 23: 14:00-18:00, interrupt 60, delta 180
$ws.Range("B23").Value = 0.58333333333333337
$ws.Range("C23").Value = 0.75
$ws.Range("D23").Value = 60
$ws.Range("E23").Value = 180

# Row 24: 20:00-22:00, interrupt 30, delta 90
$ws.Range("B24").Value = 0.83333333333333337
$ws.Range("C24").Value = 0.91666666666666663
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 90

# ---- Activity column (F), entered row 24, then 23, then 21, then 22 ----
$ws.Range("F24").Value = "Key entity class 설계"
$c = $ws.Range("F24").Characters(18, 2)
$c.Font.Name = "돋움"
$c.Font.Size = 10
$c.Font.ColorIndex = -4105

$ws.Range("F23").Value = "sequence diagram 설계"

$ws.Range("F21").Value = "node JS 강의"
$c = $ws.Range("F21").Characters(9, 2)
$c.Font.Name = "돋움"
$c.Font.Size = 10
$c.Font.ColorIndex = -4105

# This text already exists verbatim (with identical rich-text run split) as
# an earlier shared string in the workbook ("Node JS 강의"), so a plain
# assignment naturally dedupes against it rather than creating a new entry.
$ws.Range("F22").Value = "Node JS 강의"

# The author's last on-screen selection before saving was F28.
$ws.Range("F28").Select() | Out-Null
